$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.359656
$ws.Range("H2").Value = 7.078968
$ws.Range("I2").Value = 0.1135804410355361
$ws.Range("J2").Value = 0.1135804410355361
$ws.Range("O2").Value = 0.4715772180415435
$ws.Range("P2").Value = 0.4715772180415434
$ws.Range("Q2").Value = 0.515674502928
$ws.Range("R2").Value = 4.641070526352
$ws.Range("S2").Value = 0.05356194840746969
$ws.Range("T2").Value = 0.05356194840746968

# Row 3
$ws.Range("G3").Value = 2.359656
$ws.Range("H3").Value = 7.078968
$ws.Range("I3").Value = 0.1135804410355361
$ws.Range("J3").Value = 0.1135804410355361
$ws.Range("M3").Value = 0.2448813333333333
$ws.Range("N3").Value = 0.734644
$ws.Range("O3").Value = 0.5284227819584566
$ws.Range("P3").Value = 0.5284227819584566
$ws.Range("Q3").Value = 0.5778357074879998
$ws.Range("R3").Value = 5.200521367392
$ws.Range("S3").Value = 0.06001849262806644
$ws.Range("T3").Value = 0.06001849262806644

# Row 4
$ws.Range("I4").Value = 0.07630393871923234
$ws.Range("J4").Value = 0.07630393871923234
$ws.Range("O4").Value = 0.4715772180415435
$ws.Range("P4").Value = 0.4715772180415434
$ws.Range("S4").Value = 0.035983199146828
$ws.Range("T4").Value = 0.03598319914682799

# Row 5
$ws.Range("I5").Value = 0.07630393871923234
$ws.Range("J5").Value = 0.07630393871923234
$ws.Range("M5").Value = 0.2448813333333333
$ws.Range("N5").Value = 0.734644
$ws.Range("O5").Value = 0.5284227819584566
$ws.Range("P5").Value = 0.5284227819584566
$ws.Range("Q5").Value = 0.3881930727857777
$ws.Range("R5").Value = 3.493737655072
$ws.Range("S5").Value = 0.04032073957240435
$ws.Range("T5").Value = 0.04032073957240435

# Row 6
$ws.Range("G6").Value = 3.018243333333333
$ws.Range("H6").Value = 9.05473
$ws.Range("I6").Value = 0.1452810956141771
$ws.Range("J6").Value = 0.1452810956141771
$ws.Range("O6").Value = 0.4715772180415435
$ws.Range("P6").Value = 0.4715772180415434
$ws.Range("Q6").Value = 0.65960086158
$ws.Range("R6").Value = 5.93640775422
$ws.Range("S6").Value = 0.06851125490376112
$ws.Range("T6").Value = 0.0685112549037611

# Row 7
$ws.Range("G7").Value = 3.018243333333333
$ws.Range("H7").Value = 9.05473
$ws.Range("I7").Value = 0.1452810956141771
$ws.Range("J7").Value = 0.1452810956141771
$ws.Range("M7").Value = 0.2448813333333333
$ws.Range("N7").Value = 0.734644
$ws.Range("O7").Value = 0.5284227819584566
$ws.Range("P7").Value = 0.5284227819584566
$ws.Range("Q7").Value = 0.7391114517911109
$ws.Range("R7").Value = 6.652003066119999
$ws.Range("S7").Value = 0.076769840710416
$ws.Range("T7").Value = 0.076769840710416

# Row 8
$ws.Range("G8").Value = 1.732629
$ws.Range("H8").Value = 5.197887
$ws.Range("I8").Value = 0.083398921694925
$ws.Range("J8").Value = 0.083398921694925
$ws.Range("O8").Value = 0.4715772180415435
$ws.Range("P8").Value = 0.4715772180415434
$ws.Range("Q8").Value = 0.378645276402
$ws.Range("R8").Value = 3.407807487618
$ws.Range("S8").Value = 0.03932903148055725
$ws.Range("T8").Value = 0.03932903148055725

# Row 9
$ws.Range("G9").Value = 1.732629
$ws.Range("H9").Value = 5.197887
$ws.Range("I9").Value = 0.083398921694925
$ws.Range("J9").Value = 0.083398921694925
$ws.Range("M9").Value = 0.2448813333333333
$ws.Range("N9").Value = 0.734644
$ws.Range("O9").Value = 0.5284227819584566
$ws.Range("P9").Value = 0.5284227819584566
$ws.Range("Q9").Value = 0.4242884996919999
$ws.Range("R9").Value = 3.818596497228
$ws.Range("S9").Value = 0.04406989021436775
$ws.Range("T9").Value = 0.04406989021436775

# Row 10
$ws.Range("G10").Value = 9.589644
$ws.Range("H10").Value = 28.768932
$ws.Range("I10").Value = 0.4615910094072114
$ws.Range("J10").Value = 0.4615910094072114
$ws.Range("O10").Value = 0.4715772180415435
$ws.Range("P10").Value = 0.4715772180415434
$ws.Range("Q10").Value = 2.095701620472
$ws.Range("R10").Value = 18.861314584248
$ws.Range("S10").Value = 0.2176758040892406
$ws.Range("T10").Value = 0.2176758040892406

# Row 11
$ws.Range("G11").Value = 9.589644
$ws.Range("H11").Value = 28.768932
$ws.Range("I11").Value = 0.4615910094072114
$ws.Range("J11").Value = 0.4615910094072114
$ws.Range("M11").Value = 0.2448813333333333
$ws.Range("N11").Value = 0.734644
$ws.Range("O11").Value = 0.5284227819584566
$ws.Range("P11").Value = 0.5284227819584566
$ws.Range("Q11").Value = 2.348324808912
$ws.Range("R11").Value = 21.134923280208
$ws.Range("S11").Value = 0.2439152053179708
$ws.Range("T11").Value = 0.2439152053179708

# Row 12
$ws.Range("G12").Value = 2.489795
$ws.Range("H12").Value = 7.469385
$ws.Range("I12").Value = 0.1198445935289181
$ws.Range("J12").Value = 0.1198445935289181
$ws.Range("O12").Value = 0.4715772180415435
$ws.Range("P12").Value = 0.4715772180415434
$ws.Range("Q12").Value = 0.54411481971
$ws.Range("R12").Value = 4.897033377390001
$ws.Range("S12").Value = 0.05651598001368674
$ws.Range("T12").Value = 0.05651598001368673

# Row 13
$ws.Range("G13").Value = 2.489795
$ws.Range("H13").Value = 7.469385
$ws.Range("I13").Value = 0.1198445935289181
$ws.Range("J13").Value = 0.1198445935289181
$ws.Range("M13").Value = 0.2448813333333333
$ws.Range("N13").Value = 0.734644
$ws.Range("O13").Value = 0.5284227819584566
$ws.Range("P13").Value = 0.5284227819584566
$ws.Range("Q13").Value = 0.6097043193266666
$ws.Range("R13").Value = 5.48733887394
$ws.Range("S13").Value = 0.06332861351523132
$ws.Range("T13").Value = 0.06332861351523132
